$wb = $excel.ActiveWorkbook

# Update "想去人数" (F column) values across all four sheets to match the
# regenerated site data (gh-pages output at commit 456a3b4).

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 236
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 662
$ws.Range("F8").Value = 12620
$ws.Range("F9").Value = 13025
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("F12").Value = 5503
$ws.Range("F14").Value = 109
$ws.Range("F16").Value = 195
$ws.Range("F17").Value = 1439
$ws.Range("F20").Value = 0
$ws.Range("F24").Value = 511
$ws.Range("F27").Value = 261
$ws.Range("F28").Value = 0
$ws.Range("F32").Value = 1007
$ws.Range("F33").Value = 160
$ws.Range("F34").Value = 51
$ws.Range("F35").Value = 105
$ws.Range("F36").Value = 0
$ws.Range("F37").Value = 4451
$ws.Range("F38").Value = 276
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 605
$ws.Range("F42").Value = 2052
$ws.Range("F43").Value = 35
$ws.Range("F44").Value = 915
$ws.Range("F45").Value = 300
$ws.Range("F46").Value = 0
$ws.Range("F48").Value = 4303
$ws.Range("F49").Value = 0

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 15
$ws.Range("F4").Value = 36
$ws.Range("F5").Value = 12
$ws.Range("F7").Value = 40
$ws.Range("F12").Value = 0
$ws.Range("F17").Value = 6
$ws.Range("F18").Value = 11
$ws.Range("F19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("F21").Value = 13
$ws.Range("F22").Value = 77
$ws.Range("F23").Value = 76
$ws.Range("F25").Value = 108
$ws.Range("F26").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("F30").Value = 0

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 0

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 97
$ws.Range("F8").Value = 12620
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 1318
$ws.Range("F12").Value = 1297
$ws.Range("F13").Value = 5504
$ws.Range("F14").Value = 916
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 1439
$ws.Range("F18").Value = 364
$ws.Range("F20").Value = 1042
$ws.Range("F23").Value = 0
$ws.Range("F24").Value = 3032
$ws.Range("F25").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F31").Value = 14
$ws.Range("F32").Value = 0
$ws.Range("F33").Value = 160
$ws.Range("F34").Value = 51
$ws.Range("F35").Value = 105
$ws.Range("F37").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("F39").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 646
$ws.Range("F43").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 300
$ws.Range("F46").Value = 44
$ws.Range("F47").Value = 33
$ws.Range("F48").Value = 4303
